$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 162-171 (date / price column changes) ---
$ws.Range("D162").Value = 44585
$ws.Range("K162").Value = 2000
$ws.Range("L162").Value = 2000
$ws.Range("M162").Value = 2000
$ws.Range("P162").Value = 2000

$ws.Range("D163").Value = 44585
$ws.Range("K163").Value = 1500
$ws.Range("L163").Value = 1500
$ws.Range("M163").Value = 1500
$ws.Range("P163").Value = 1500

$ws.Range("D164").Value = 44585
$ws.Range("K164").Value = 1000
$ws.Range("L164").Value = 1000
$ws.Range("M164").Value = 1000
$ws.Range("P164").Value = 1000

$ws.Range("D165").Value = 44560
$ws.Range("J165").Value = 2000
$ws.Range("K165").Value = 2500
$ws.Range("L165").Value = 2500
$ws.Range("M165").Value = 2500
$ws.Range("P165").Value = 2500

$ws.Range("D166").Value = 44560
$ws.Range("J166").Value = 3000
$ws.Range("K166").Value = 2000
$ws.Range("L166").Value = 2000
$ws.Range("M166").Value = 2000
$ws.Range("P166").Value = 2000

$ws.Range("D167").Value = 44560
$ws.Range("J167").Value = 2000
$ws.Range("K167").Value = 1500
$ws.Range("L167").Value = 1500
$ws.Range("M167").Value = 1500
$ws.Range("P167").Value = 1500

$ws.Range("D168").Value = 44189
$ws.Range("J168").Value = 1000
$ws.Range("K168").Value = 3500
$ws.Range("L168").Value = 3500
$ws.Range("M168").Value = 3500
$ws.Range("P168").Value = 3500

$ws.Range("D169").Value = 44189
$ws.Range("J169").Value = 2500
$ws.Range("K169").Value = 3000
$ws.Range("L169").Value = 3000
$ws.Range("M169").Value = 3000
$ws.Range("P169").Value = 3000

$ws.Range("D170").Value = 44189
$ws.Range("J170").Value = 1500
$ws.Range("K170").Value = 2500
$ws.Range("L170").Value = 2500
$ws.Range("M170").Value = 2500
$ws.Range("P170").Value = 2500

$ws.Range("D171").Value = 44209
$ws.Range("K171").Value = 2000
$ws.Range("L171").Value = 2000
$ws.Range("M171").Value = 2000
$ws.Range("P171").Value = 2000

# --- Append new rows 172-175 ---
# Row 172
$ws.Range("A172").Value = 5
$ws.Range("B172").Value = "Macroferia Regional de Talca"
$ws.Range("C172").Value = "Maule"
$ws.Range("D172").Value = 44209
$ws.Range("E172").Value = 7
$ws.Range("F172").Value = 100112028
$ws.Range("G172").Value = "Sandia"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 5000
$ws.Range("K172").Value = 1600
$ws.Range("L172").Value = 1600
$ws.Range("M172").Value = 1600
$ws.Range("N172").Value = "$/unidad"
$ws.Range("O172").Value = "Región del Maule"
$ws.Range("P172").Value = 1600
$ws.Range("Q172").Value = 1
$ws.Range("R172").Value = "Hortaliza"
$ws.Range("D172").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 173
$ws.Range("A173").Value = 5
$ws.Range("B173").Value = "Macroferia Regional de Talca"
$ws.Range("C173").Value = "Maule"
$ws.Range("D173").Value = 44209
$ws.Range("E173").Value = 7
$ws.Range("F173").Value = 100112028
$ws.Range("G173").Value = "Sandia"
$ws.Range("H173").Value = "Sin especificar"
$ws.Range("I173").Value = "Segunda"
$ws.Range("J173").Value = 5000
$ws.Range("K173").Value = 1200
$ws.Range("L173").Value = 1200
$ws.Range("M173").Value = 1200
$ws.Range("N173").Value = "$/unidad"
$ws.Range("O173").Value = "Región del Maule"
$ws.Range("P173").Value = 1200
$ws.Range("Q173").Value = 1
$ws.Range("R173").Value = "Hortaliza"
$ws.Range("D173").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 174
$ws.Range("A174").Value = 5
$ws.Range("B174").Value = "Macroferia Regional de Talca"
$ws.Range("C174").Value = "Maule"
$ws.Range("D174").Value = 44554
$ws.Range("E174").Value = 7
$ws.Range("F174").Value = 100112028
$ws.Range("G174").Value = "Sandia"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Extra"
$ws.Range("J174").Value = 3000
$ws.Range("K174").Value = 2500
$ws.Range("L174").Value = 2500
$ws.Range("M174").Value = 2500
$ws.Range("N174").Value = "$/unidad"
$ws.Range("O174").Value = "Región del Maule"
$ws.Range("P174").Value = 2500
$ws.Range("Q174").Value = 1
$ws.Range("R174").Value = "Hortaliza"
$ws.Range("D174").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 175
$ws.Range("A175").Value = 5
$ws.Range("B175").Value = "Macroferia Regional de Talca"
$ws.Range("C175").Value = "Maule"
$ws.Range("D175").Value = 44554
$ws.Range("E175").Value = 7
$ws.Range("F175").Value = 100112028
$ws.Range("G175").Value = "Sandia"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 3000
$ws.Range("K175").Value = 1900
$ws.Range("L175").Value = 1900
$ws.Range("M175").Value = 1900
$ws.Range("N175").Value = "$/unidad"
$ws.Range("O175").Value = "Región del Maule"
$ws.Range("P175").Value = 1900
$ws.Range("Q175").Value = 1
$ws.Range("R175").Value = "Hortaliza"
$ws.Range("D175").NumberFormat = "YYYY-MM-DD HH:MM:SS"
